# 5.2.1 Data.xlsx edit: dynamic-method table gets two new computed columns
# (a model/theoretical Ik column and the residual Ik - model column), the
# column B header typo "Iк, мА" is fixed to "Iк, мкА", two new header labels
# "I0, мкА" / "I~, мкА" are introduced, and a small summary table is added
# below the main one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Row 19 headers: fix the "Iк, мА" typo (B19) and add the two new
#    header columns (C/D, G/H, K/L) for each of the three side-by-side
#    mini tables.
# ---------------------------------------------------------------------
$ws.Range("B19").Value = "Iк, мкА"
$ws.Range("C19").Value = "I0, мкА"
$ws.Range("D19").Value = "I~, мкА"

$ws.Range("F19").Value = "Iк, мкА"
$ws.Range("G19").Value = "I0, мкА"
$ws.Range("H19").Value = "I~, мкА"

$ws.Range("J19").Value = "Iк, мкА"
$ws.Range("K19").Value = "I0, мкА"
$ws.Range("L19").Value = "I~, мкА"

# ---------------------------------------------------------------------
# 2) Column B (rows 20:65) used to store the current in amps (e.g.
#    0.028); retype as microamps (e.g. 28) to match the corrected header.
# ---------------------------------------------------------------------
$bvals = @{20=28;21=32;22=38;23=53;24=60;25=68;26=81;27=91;28=107;29=118;30=128;31=137;32=141;33=142;34=140;35=137;36=131;37=123;38=77;39=65;40=67;41=87;42=108;43=121;44=131;45=145;46=153;47=168;48=171;49=187;50=200;51=202;52=201;53=201;54=187;55=174;56=169;57=163;58=175;59=192;60=219;61=231;62=237;63=240;64=248;65=261}

for ($r = 20; $r -le 65; $r++) {
    $ws.Range("B$r").Value = $bvals[$r]
}

# ---------------------------------------------------------------------
# 3) Column D (new): residual = measured - model. Fill before column C so
#    this shared-formula group is allocated first (matches the original
#    file's existing shared group already occupying column C).
# ---------------------------------------------------------------------
$ws.Range("D20").Formula = "=B20-C20"
$ws.Range("D21:D65").Formula = "=B21-C21"

# ---------------------------------------------------------------------
# 4) Column C: replace the old "=B*1000" formula with the fitted
#    polynomial model of Ik vs Va.
# ---------------------------------------------------------------------
$ws.Range("C20").Formula = "=12.5321106194744 + 6.40918573732954 * A20 - 0.105645947480177 * A20*A20 + 0.00172074354709666*A20*A20*A20  - 1.22543389562474 * A20*A20*A20*A20 / 100000"
$ws.Range("C21:C65").Formula = "=12.5321106194744 + 6.40918573732954 * A21 - 0.105645947480177 * A21*A21 + 0.00172074354709666*A21*A21*A21  - 1.22543389562474 * A21*A21*A21*A21 / 100000"

# ---------------------------------------------------------------------
# 5) New summary table appended below (rows 69-79).
# ---------------------------------------------------------------------
$ws.Range("A69").Value = "Va, В"
$ws.Range("B69").Value = "I0, мкА"
$ws.Range("E69").Value = "Va, В"
$ws.Range("F69").Value = "I0, мкА"
$ws.Range("I69").Value = "Va, В"
$ws.Range("J69").Value = "I0, мкА"

$ws.Range("A70").Value = 2.95
$ws.Range("B70").Value = 28
$ws.Range("E70").Value = 2.35
$ws.Range("F70").Value = 13

$ws.Range("A71").Value = 3.65
$ws.Range("B71").Value = 32
$ws.Range("E71").Value = 4.22
$ws.Range("F71").Value = 22

$ws.Range("A72").Value = 4.55
$ws.Range("B72").Value = 38
$ws.Range("E72").Value = 5.76
$ws.Range("F72").Value = 32

$ws.Range("A73").Value = 6.33
$ws.Range("B73").Value = 53
$ws.Range("E73").Value = 6.44
$ws.Range("F73").Value = 36

$ws.Range("A74").Value = 7.18
$ws.Range("B74").Value = 60
$ws.Range("E74").Value = 23.3
$ws.Range("F74").Value = 88

$ws.Range("A75").Value = 22.5
$ws.Range("B75").Value = 100
$ws.Range("E75").Value = 30.8
$ws.Range("F75").Value = 110

$ws.Range("A76").Value = 30.66
$ws.Range("B76").Value = 145

$ws.Range("A77").Value = 40.18
$ws.Range("B77").Value = 187

$ws.Range("A78").Value = 56.51
$ws.Range("B78").Value = 219

$ws.Range("A79").Value = 68
$ws.Range("B79").Value = 240

# ---------------------------------------------------------------------
# 6) View state: scroll position + active selection moved down to the
#    newly added rows.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("E74:F74").Select()
